{"js": "// Replace the 25 \"a\u00f7b=c, r\" answer strings inside the table with their\n// updated values, as described by the diff. Every table cell in the\n// document gets a brand new value, so we do a direct exact-text\n// search-and-replace for each (old, new) pair. Because every \"old\"\n// value is unique in the document, this is unambiguous.\nconst replacements = [\n  [\"995\u00f72=497, 1\", \"482\u00f72=241, 0\"],\n  [\"354\u00f72=177, 0\", \"964\u00f74=241, 0\"],\n  [\"271\u00f72=135, 1\", \"685\u00f79=76, 1\"],\n  [\"956\u00f78=119, 4\", \"195\u00f78=24, 3\"],\n  [\"371\u00f74=92, 3\", \"773\u00f79=85, 8\"],\n  [\"671\u00f76=111, 5\", \"895\u00f72=447, 1\"],\n  [\"901\u00f76=150, 1\", \"789\u00f76=131, 3\"],\n  [\"496\u00f78=62, 0\", \"482\u00f72=241, 0\"],\n  [\"606\u00f76=101, 0\", \"869\u00f79=96, 5\"],\n  [\"752\u00f77=107, 3\", \"897\u00f76=149, 3\"],\n  [\"796\u00f72=398, 0\", \"207\u00f74=51, 3\"],\n  [\"753\u00f75=150, 3\", \"564\u00f77=80, 4\"],\n  [\"896\u00f77=128, 0\", \"158\u00f75=31, 3\"],\n  [\"691\u00f74=172, 3\", \"197\u00f72=98, 1\"],\n  [\"123\u00f75=24, 3\", \"377\u00f72=188, 1\"],\n  [\"699\u00f73=233, 0\", \"237\u00f73=79, 0\"],\n  [\"770\u00f76=128, 2\", \"560\u00f74=140, 0\"],\n  [\"615\u00f72=307, 1\", \"809\u00f78=101, 1\"],\n  [\"331\u00f78=41, 3\", \"882\u00f72=441, 0\"],\n  [\"408\u00f74=102, 0\", \"792\u00f76=132, 0\"],\n  [\"169\u00f78=21, 1\", \"904\u00f73=301, 1\"],\n  [\"373\u00f74=93, 1\", \"293\u00f73=97, 2\"],\n  [\"287\u00f75=57, 2\", \"542\u00f76=90, 2\"],\n  [\"148\u00f76=24, 4\", \"514\u00f76=85, 4\"],\n  [\"392\u00f73=130, 2\", \"453\u00f76=75, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('995\u00f72=497, 1', '482\u00f72=241, 0'),\n    @('354\u00f72=177, 0', '964\u00f74=241, 0'),\n    @('271\u00f72=135, 1', '685\u00f79=76, 1'),\n    @('956\u00f78=119, 4', '195\u00f78=24, 3'),\n    @('371\u00f74=92, 3', '773\u00f79=85, 8'),\n    @('671\u00f76=111, 5', '895\u00f72=447, 1'),\n    @('901\u00f76=150, 1', '789\u00f76=131, 3'),\n    @('496\u00f78=62, 0', '482\u00f72=241, 0'),\n    @('606\u00f76=101, 0', '869\u00f79=96, 5'),\n    @('752\u00f77=107, 3', '897\u00f76=149, 3'),\n    @('796\u00f72=398, 0', '207\u00f74=51, 3'),\n    @('753\u00f75=150, 3', '564\u00f77=80, 4'),\n    @('896\u00f77=128, 0', '158\u00f75=31, 3'),\n    @('691\u00f74=172, 3', '197\u00f72=98, 1'),\n    @('123\u00f75=24, 3', '377\u00f72=188, 1'),\n    @('699\u00f73=233, 0', '237\u00f73=79, 0'),\n    @('770\u00f76=128, 2', '560\u00f74=140, 0'),\n    @('615\u00f72=307, 1', '809\u00f78=101, 1'),\n    @('331\u00f78=41, 3', '882\u00f72=441, 0'),\n    @('408\u00f74=102, 0', '792\u00f76=132, 0'),\n    @('169\u00f78=21, 1', '904\u00f73=301, 1'),\n    @('373\u00f74=93, 1', '293\u00f73=97, 2'),\n    @('287\u00f75=57, 2', '542\u00f76=90, 2'),\n    @('148\u00f76=24, 4', '514\u00f76=85, 4'),\n    @('392\u00f73=130, 2', '453\u00f76=75, 3'),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $oldText,  # FindText\n        $true,     # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap = wdFindContinue\n        $false,    # Format\n        $newText,  # ReplaceWith\n        2          # Replace = wdReplaceAll\n    )\n}\n"}
